$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Board")

# Fix S6: was referencing N7 instead of N6
$ws.Range("S6").Formula = "=L6&M6&N6&O6&P6&Q6"

# Fix S7: was referencing N8 instead of N7
$ws.Range("S7").Formula = "=L7&M7&N7&O7&P7&Q7"

# Fix S8: was referencing broken #REF! instead of N8
$ws.Range("S8").Formula = "=L8&M8&N8&O8&P8&Q8"

# Update the sheet view: scroll position and selection
$ws.Range("G14").Select()

$excel.ActiveWindow.ScrollRow = 2
